# "Commit after Excel data updated"
#
# Sheet1 holds username/password rows. The username in column A changes
# from "alagappan.n@vmokshgroup.com" to "admin@mydomain.com" for every
# data row (A2:A4), and the A3/A4 hyperlinks are removed (only A2 keeps a
# mailto: link, matching the new layout). Column B (the passwords) is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The COM shim only supports bulk-deleting every hyperlink on the sheet
# (Range.Hyperlinks.Delete()/single Hyperlink.Delete() are no-ops here), so
# clear them all and rebuild only the ones that should remain.
$ws.Hyperlinks.Delete()

# Column A: replace the old username with the new one on every data row.
$ws.Range("A2").Value = "admin@mydomain.com"
$ws.Range("A3").Value = "admin@mydomain.com"
$ws.Range("A4").Value = "admin@mydomain.com"

# Re-add the hyperlinks that remain after the update: A2 (new username) and
# all three passwords in column B (unchanged targets).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@mydomain.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Power@123")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Power@1234")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Powe@1243")

# Hyperlinks.Add() stamps a freshly-minted style on the anchor cell; put the
# cells back on the shared "Hyperlink" style so formatting matches the rest
# of the sheet.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
